$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update L6: CV -> VC
$ws.Range("L6").Value = "VC"
$ws.Range("M6").Value = 4
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0

# Update L7: VC -> V
$ws.Range("L7").Value = "V"
$ws.Range("M7").Value = 2
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0

# New row 8: CV
$ws.Range("L8").Value = "CV"
$ws.Range("M8").Value = 4
$ws.Range("N8").Value = 7
$ws.Range("O8").Value = 5

# New row 10: CVV (set before row 9's CVC so the shared-string table
# gets new unique strings in the order CVV, CVC, CCV - matching the target)
$ws.Range("L10").Value = "CVV"
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 4

# New row 9: CVC
$ws.Range("L9").Value = "CVC"
$ws.Range("M9").Value = 4
$ws.Range("N9").Value = 5
$ws.Range("O9").Value = 2

# New row 11: CCV
$ws.Range("L11").Value = "CCV"
$ws.Range("M11").Value = 1
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 2

# Re-affirm H28/I28 totals formulas (column totals row)
$ws.Range("H28").Formula = "=SUM(H2:H27)"
$ws.Range("I28").Formula = "=SUM(I2:I27)"

# Update sheet view: change selection (also resets any scrolled topLeftCell)
$ws.Range("L14").Select()
